# Updates cryptos list values (price + 1h volume change) per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.958.80"
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").Value = "'1.884.85"
$ws.Range("E3").Value = "  +1.34%  "
$ws.Range("D4").Value = "'0.9998"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'305.33"
$ws.Range("E5").Value = "  +0.28%  "
$ws.Range("D6").Value = "'0.9994"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").Value = "'0.5151"
$ws.Range("E7").Value = "  +2.45%  "
$ws.Range("D8").Value = "'0.3736"
$ws.Range("E8").Value = "  +2.50%  "
$ws.Range("D9").Value = "'0.07185"
$ws.Range("E9").Value = "  +0.22%  "
$ws.Range("D10").Value = "'21.04"
$ws.Range("E10").Value = "  +1.69%  "
$ws.Range("E11").Value = "  +0.62%  "
$ws.Range("D12").Value = "'0.07646"
$ws.Range("E12").Value = "  +2.18%  "
$ws.Range("D13").Value = "'1.864.57"
$ws.Range("E13").Value = "  +0.06%  "
$ws.Range("D14").Value = "'93.68"
$ws.Range("E14").Value = "  +0.10%  "
$ws.Range("D15").Value = "'5.228"
$ws.Range("E15").Value = "  +0.01%  "
$ws.Range("D16").Value = "'1.000"
$ws.Range("E16").Value = "  -0.02%  "
$ws.Range("D17").Value = "'0.000008471"
$ws.Range("E17").Value = "  -0.23%  "
$ws.Range("E18").Value = "  +1.45%  "
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("D20").Value = "'26.996.18"
$ws.Range("E20").Value = "  +0.23%  "
$ws.Range("D21").Value = "'5.033"
$ws.Range("E21").Value = "  +0.18%  "
$ws.Range("D22").Value = "'2.113.15"
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").Value = "'10.54"
$ws.Range("E23").Value = "  +1.53%  "
$ws.Range("D24").Value = "'6.375"
$ws.Range("E24").Value = "  -0.58%  "
$ws.Range("D25").Value = "'2.286"
$ws.Range("E25").Value = "  +9.74%  "
$ws.Range("D26").Value = "'146.14"
$ws.Range("E26").Value = "  -0.95%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'18.00"
$ws.Range("E27").Value = "  +0.82%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "'1.723"
$ws.Range("E28").Value = "  -3.32%  "
$ws.Range("E29").Value = "  +0.73%  "
$ws.Range("D30").Value = "'4.904"
$ws.Range("E30").Value = "  +5.08%  "
$ws.Range("E31").Value = "  +1.80%  "
$ws.Range("D32").Value = "'0.09176"
$ws.Range("E32").Value = "  -0.45%  "
$ws.Range("E33").Value = "  -2.11%  "
$ws.Range("D34").Value = "'1.227"
$ws.Range("E34").Value = "  +6.75%  "
$ws.Range("D35").Value = "'0.7650"
$ws.Range("E35").Value = "  +2.45%  "
$ws.Range("D36").Value = "'2.985"
$ws.Range("E36").Value = "  +0.65%  "
$ws.Range("D37").Value = "'3.259"
$ws.Range("E37").Value = "  +0.23%  "
$ws.Range("E38").Value = "  +0.61%  "
$ws.Range("D39").Value = "'0.5562"
$ws.Range("E39").Value = "  +0.16%  "
$ws.Range("D40").Value = "'0.01983"
$ws.Range("E40").Value = "  -0.97%  "
$ws.Range("D41").Value = "'1.070"
$ws.Range("E41").Value = "  +0.20%  "
$ws.Range("D42").Value = "'8.998"
$ws.Range("E42").Value = "  +5.73%  "
$ws.Range("D43").Value = "'6.594"
$ws.Range("E43").Value = "  +0.79%  "
$ws.Range("D44").Value = "'118.38"
$ws.Range("E44").Value = "  +0.88%  "
$ws.Range("E45").Value = "  +1.84%  "
$ws.Range("D46").Value = "'0.4801"
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").Value = "'0.9991"
$ws.Range("E47").Value = "  -0.03%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'10.09"
$ws.Range("E48").Value = "  +0.88%  "
$ws.Range("E49").Value = "  +2.44%  "
$ws.Range("D50").Value = "'37.58"
$ws.Range("E50").Value = "  +2.55%  "
$ws.Range("D51").Value = "'63.72"
$ws.Range("E51").Value = "  +1.24%  "
